$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "2024-02-20 Tuesday" "2024-02-21 Wednesday"

Replace-Text "14×70=" "31×48="
Replace-Text "59×28=" "46×39="
Replace-Text "68×81=" "36×57="
Replace-Text "71×36=" "79×18="
Replace-Text "58×80=" "48×16="

Replace-Text "79×54=" "80×90="
Replace-Text "43×85=" "52×43="
Replace-Text "33×36=" "85×57="
Replace-Text "53×45=" "31×84="
Replace-Text "42×67=" "92×63="

Replace-Text "47×17=" "39×30="
Replace-Text "82×85=" "21×49="
Replace-Text "45×33=" "16×13="
Replace-Text "20×71=" "35×17="
Replace-Text "23×78=" "15×92="

Replace-Text "42×34=" "68×24="
Replace-Text "95×66=" "27×54="
Replace-Text "48×50=" "46×93="
Replace-Text "76×11=" "74×89="
Replace-Text "19×84=" "98×47="

Replace-Text "59×54=" "87×33="
Replace-Text "68×83=" "13×62="
Replace-Text "68×65=" "81×78="
Replace-Text "36×91=" "30×88="
Replace-Text "29×82=" "95×71="
